$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "TEST_STANDARD_ROW"
$wsStandard = $wb.Worksheets.Item("Sheet1")
$wsStandard.Name = "TEST_STANDARD_ROW"

# Add the "YO" named range, pointing at column A of the renamed sheet
$wb.Names.Add("YO", '=TEST_STANDARD_ROW!$A:$A')

# Add new data to the TEST_SCALAR_INPUT sheet (kwh / 100 row, plus dollar value)
$wsScalar = $wb.Worksheets.Item("TEST_SCALAR_INPUT")
$wsScalar.Range("B2").Value = 1
$wsScalar.Range("A3").Value = "kwh"
$wsScalar.Range("B3").Value = 100

# Update selections on each sheet (order matters: last activated sheet keeps tabSelected)
$wsTest = $wb.Worksheets.Item("TEST_SHEET")
[void]$wsTest.Range("E12").Select()

$wsTracker = $wb.Worksheets.Item("TEXEL_SHEET_TRACKER")
[void]$wsTracker.Range("G37").Select()

[void]$wsScalar.Range("I19").Select()

[void]$wsStandard.Range("K19").Select()
